{"js": "// Apply the Progress-Report revision described by the commit diff.\n// Strategy: locate each target paragraph by its current (pre-edit) text,\n// then rewrite its content with the post-edit wording. Finally remove the\n// trailing empty paragraph that the author deleted.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findParagraph(startsWith) {\n  for (const p of items) {\n    if (p.text.indexOf(startsWith) === 0) return p;\n  }\n  throw new Error(\"Paragraph not found: \" + startsWith);\n}\n\n// --- Change 1: \"Activities versus the Schedule\" paragraph -----------------\nconst p1 = findParagraph(\n  \"I am on the schedule from the beginning of the Semester\"\n);\np1.getRange(\"Whole\").insertText(\n  \"I am right on schedule from the beginning of the Semester to finish my \" +\n    \"project and present for next week. I have a waterproof (DS18B20) \" +\n    \"Temperature Sensor.\",\n  \"Replace\"\n);\n\n// --- Change 2: \"Current Progress\" paragraph --------------------------------\nconst p2 = findParagraph(\"My Sensor is working very well\");\np2.getRange(\"Whole\").insertText(\n  \"My Sensor is working very well and it shows the readings of temperature \" +\n    \"in Celsius and Fahrenheit on Computer Screen. I also tested my sensor \" +\n    \"under water and in direct flame to test the manufacturers claim and \" +\n    \"it passed with flying colours.\",\n  \"Replace\"\n);\n\n// --- Change 3: \"Problems and Opportunities\" paragraph ----------------------\nconst p3 = findParagraph(\"Before Christmas Holidays\");\np3.getRange(\"Whole\").insertText(\n  \"One of the main problem that I had faced is before Christmas holidays, \" +\n    \"I thought I should get a 3d printed case for my sensor and raspberry \" +\n    \"pi so that it looks nice while I am showing demo in class. But \" +\n    \"Anthony printed the wrong case and I got my correct case on 2nd \" +\n    \"January. The Second problem which I am facing right now is that I am \" +\n    \"trying to add an LCD screen as an add-on to my sensor but I am not \" +\n    \"able to figure out why it is not working. Currently, I have stopped \" +\n    \"working on my project because of other courses but hopefully my LCD \" +\n    \"screen will also be up and running by the start of next semester.  \",\n  \"Replace\"\n);\n\n// --- Change 4: \"Financial Status\" paragraph --------------------------------\n// Keep the existing sentence (and the _GoBack bookmark that trails it)\n// intact, just add the new lead-in and trailing sentences around it.\nconst p4 = findParagraph(\"I didn\\u2019t buy any extra stuff\");\nconst p4Range = p4.getRange(\"Whole\");\np4Range.insertText(\n  \"Financially, I am on track to meet the budget that I set at start of \" +\n    \"the semester. \",\n  \"Start\"\n);\np4.getRange(\"Whole\").insertText(\n  \" The only thing for which I paid more was my sensor because when I \" +\n    \"ordered my sensor on Amazon they were only delivering pack of three. \" +\n    \"So, at that time I didn\\u2019t have any option so I bought three \" +\n    \"temperature sensors instead of one.\",\n  \"End\"\n);\n\nawait context.sync();\n\n// --- Change 5: drop the trailing blank paragraph ---------------------------\nconst trailing = body.paragraphs;\ntrailing.load(\"items/text\");\nawait context.sync();\nconst tItems = trailing.items;\nconst lastPara = tItems[tItems.length - 1];\nif (lastPara.text.trim() === \"\") {\n  lastPara.delete();\n}\n\nawait context.sync();\n", "ps1": "# Apply the Progress-Report revision described by the commit diff.\n# Strategy: use Find/Replace against the exact pre-edit wording of each\n# paragraph so the edit is anchored on content rather than a fragile\n# paragraph index, then drop the trailing blank paragraph the author\n# removed.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\n# --- Change 1: \"Activities versus the Schedule\" paragraph -----------------\nReplace-Text `\n    \"I am on the schedule from the beginning of the Semester and I am ready for the demo of my Temperature Sensor (DS18B20).\" `\n    \"I am right on schedule from the beginning of the Semester to finish my project and present for next week. I have a waterproof (DS18B20) Temperature Sensor.\"\n\n# --- Change 2: \"Current Progress\" paragraph --------------------------------\nReplace-Text `\n    \"My Sensor is working very well and it shows the readings of temperature in Celsius and Fahrenheit. I also tested my sensor under water and in direct flame to test the manufacturers claim and it passed with flying \" `\n    \"My Sensor is working very well and it shows the readings of temperature in Celsius and Fahrenheit on Computer Screen. I also tested my sensor under water and in direct flame to test the manufacturers claim and it passed with flying \"\n\n# --- Change 3: \"Problems and Opportunities\" paragraph ----------------------\nReplace-Text `\n    \"Before Christmas Holidays, I thought I should get a 3d printed case for my sensor and raspberry pi so that it looks nice while I am showing demo in class. But Anthony printed the wrong case and I got my correct case on 2nd January. I was trying to add an LCD screen as an add-on to my sensor but i am not able to figure out why it is not working. Currently, I have stopped working on my project because of other courses but hopefully my LCD screen will also be up and running by the start of next semester.  \" `\n    \"One of the main problem that I had faced is before Christmas holidays, I thought I should get a 3d printed case for my sensor and raspberry pi so that it looks nice while I am showing demo in class. But Anthony printed the wrong case and I got my correct case on 2nd January. The Second problem which I am facing right now is that I am trying to add an LCD screen as an add-on to my sensor but I am not able to figure out why it is not working. Currently, I have stopped working on my project because of other courses but hopefully my LCD screen will also be up and running by the start of next semester.  \"\n\n# --- Change 4: \"Financial Status\" paragraph --------------------------------\n# Leave the original sentence (and the trailing _GoBack bookmark) in place;\n# just wrap it with the new lead-in / trailing sentences.\nReplace-Text `\n    \"I didn\u2019t buy any extra stuff for my project. Therefore, I am not over/under budget.\" `\n    \"Financially, I am on track to meet the budget that I set at start of the semester. I didn\u2019t buy any extra stuff for my project. Therefore, I am not over/under budget. The only thing for which I paid more was my sensor because when I ordered my sensor on Amazon they were only delivering pack of three. So, at that time I didn\u2019t have any option so I bought three temperature sensors instead of one.\"\n\n# --- Change 5: drop the trailing blank paragraph ---------------------------\n$lastParaIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($lastParaIndex)\nif ($lastPara.Range.Text.Trim() -eq \"\") {\n    $lastPara.Range.Delete()\n}\n"}
